$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57; existing rows 57-69 shift down to 58-70.
$ws.Rows.Item(57).Insert()

# Populate the new row 57 with the weekly price entry.
$ws.Cells.Item(57, 1).Value = 1
$ws.Cells.Item(57, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(57, 4).Value = 45021
$ws.Cells.Item(57, 5).Value = 15
$ws.Cells.Item(57, 6).Value = 100112028
$ws.Cells.Item(57, 7).Value = "Sandia"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Segunda"
$ws.Cells.Item(57, 10).Value = 800
$ws.Cells.Item(57, 11).Value = 370
$ws.Cells.Item(57, 12).Value = 380
$ws.Cells.Item(57, 13).Value = 375
$ws.Cells.Item(57, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(57, 15).Value = "Perú"
$ws.Cells.Item(57, 16).Value = 375
$ws.Cells.Item(57, 17).Value = 1
$ws.Cells.Item(57, 18).Value = "Hortaliza"
